# Insert a new row into the weekly price table for Camote (Vega Modelo de Temuco)
# at row 212, pushing the existing rows 212:235 down to 213:236, and fill the
# new row with this week's data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 212 (existing row 212 and below shift down by one)
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row 212 with the new weekly record
$ws.Range("A212").Value2 = 10
$ws.Range("B212").Value2 = "Vega Modelo de Temuco"
$ws.Range("C212").Value2 = "La Araucanía"
$ws.Range("D212").Value2 = 45194
$ws.Range("E212").Value2 = 9
$ws.Range("F212").Value2 = 100114002
$ws.Range("G212").Value2 = "Camote"
$ws.Range("H212").Value2 = "Sin especificar"
$ws.Range("I212").Value2 = "Primera"
$ws.Range("J212").Value2 = 80
$ws.Range("K212").Value2 = 24000
$ws.Range("L212").Value2 = 24000
$ws.Range("M212").Value2 = 24000
$ws.Range("N212").Value2 = "$/caja 18 kilos"
$ws.Range("O212").Value2 = "Perú"
$ws.Range("P212").Value2 = 1333
$ws.Range("Q212").Value2 = 18
$ws.Range("R212").Value2 = "Hortaliza"

Write-Host "Row inserted and populated. New dimension rows: $($ws.UsedRange.Rows.Count)"
